$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 - "Save", reusing the formatting of the existing
# header cells (bold, centered, thin border) by copying from G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data column H2:H5, numeric values of 0 (no special style), matching
# column G's data cells.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
